# Slide 1, shape "텍스트 개체 틀 2" (id=3): second paragraph "첫째" -> "둘째"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

$fullText = $tr.Text
$idx = $fullText.IndexOf("첫째")
if ($idx -ge 0) {
    $target = $tr.Characters($idx + 1, 2)
    $target.Text = "둘째"
}
